$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-07 20:11:35", "hatespeech", "setting3", 2, 40.78162014830423),
    @("2023-10-07 20:11:35", "hatespeech", "setting5", 2, 32.85517692643104),
    @("2023-10-07 20:11:35", "hatespeech", "setting2", 2, 33.60458283673321),
    @("2023-10-07 20:11:35", "hatespeech", "def",      2, 35.3028277213472),
    @("2023-10-07 20:11:35", "hatespeech", "setting4", 2, 37.06493320621639),
    @("2023-10-07 20:11:35", "hatespeech", "setting1", 2, 42.36895803107342)
)

$startRow = 40
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}
